$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "한국과 선진국의 격차"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/korea-lacking-math/#utm_source=rss&utm_medium=rss&utm_campaign=korea-lacking-math"

$ws.Range("D30").Value = "파이썬 Multiprocessing + tqdm 활용"
$ws.Range("E30").Value = "https://john-analyst.medium.com/%ED%8C%8C%EC%9D%B4%EC%8D%AC-multiprocessing-tqdm-%ED%99%9C%EC%9A%A9-a4bc02b8790c?source=rss-2822c56b4ae3------2"

$ws.Range("D37").Value = "[Paper Review] Exploring the Limits of the Out-of-Distribution Detection"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1918&mod=document&pageid=1"

$ws.Range("D51").Value = "[flask+jinja2] 서버에서 받은 html 요소가 html 문서에서 제대로 표현되게 하려면?"
$ws.Range("E51").Value = "https://bskyvision.com/1223"

$ws.Range("D52").Value = "서울시 미세/초미세 먼지(2/N)"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2507&utm_source=rss&utm_medium=rss&utm_campaign=%25ec%2584%259c%25ec%259a%25b8%25ec%258b%259c-%25eb%25af%25b8%25ec%2584%25b8-%25ec%25b4%2588%25eb%25af%25b8%25ec%2584%25b8-%25eb%25a8%25bc%25ec%25a7%25802-n"
